$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 92; $r++) {
    $ws.Range("Z${r}").Formula = "=IF(MAX(B${r}:Y${r}) - MIN(B${r}:Y${r}) > 17, 1, 0)"
}

$ws.Range("AA2").Formula = "=SUM(Z:Z)"
$ws.Range("AA2").Style = "Good"

[void]$ws.Range("AA2").Select()
